# Subject/seed performance results for CTNet.
# Re-generated metrics: per-patient scores rounded to 2 decimals and the
# per-row/per-column "Average" columns recomputed to match, per the latest
# create_excel run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number -> [Patient1..Patient9, Average] values
$results = @{
    2 = @(0.77, 0.71, 0.76, 0.97, 0.91, 0.87, 0.91, 0.94, 0.88, 0.8577777777777778)
    3 = @(0.81, 0.7, 0.77, 0.9, 0.98, 0.86, 0.91, 0.94, 0.9, 0.8633333333333333)
    4 = @(0.79, 0.63, 0.78, 0.88, 0.97, 0.93, 0.91, 0.95, 0.89, 0.8588888888888888)
    5 = @(0.79, 0.71, 0.77, 0.9, 0.98, 0.88, 0.9, 0.94, 0.88, 0.861111111111111)
    6 = @(0.75, 0.69, 0.74, 0.86, 0.97, 0.87, 0.92, 0.94, 0.89, 0.8477777777777777)
    7 = @(0.76, 0.69, 0.79, 0.93, 0.99, 0.87, 0.9, 0.95, 0.91, 0.8655555555555555)
    8 = @(0.81, 0.71, 0.8, 0.93, 0.97, 0.86, 0.9, 0.94, 0.9, 0.8688888888888889)
    9 = @(0.77, 0.72, 0.72, 0.96, 0.97, 0.88, 0.91, 0.95, 0.89, 0.8633333333333333)
    10 = @(0.78, 0.69, 0.78, 0.9, 0.98, 0.87, 0.92, 0.93, 0.88, 0.8588888888888889)
    11 = @(0.76, 0.67, 0.77, 0.97, 0.97, 0.85, 0.91, 0.95, 0.9, 0.8611111111111112)
    12 = @(0.77, 0.7, 0.71, 0.9, 0.98, 0.87, 0.89, 0.94, 0.9, 0.8511111111111112)
    13 = @(0.73, 0.69, 0.75, 0.97, 0.96, 0.9, 0.89, 0.95, 0.87, 0.8566666666666667)
    14 = @(0.76, 0.62, 0.76, 0.91, 0.96, 0.89, 0.91, 0.94, 0.89, 0.8488888888888888)
    15 = @(0.8, 0.69, 0.75, 0.96, 0.96, 0.9, 0.92, 0.94, 0.88, 0.8666666666666667)
    16 = @(0.775, 0.687, 0.761, 0.924, 0.968, 0.879, 0.907, 0.943, 0.89, 0.8593333333333333)
}

foreach ($row in $results.Keys) {
    $rowValues = $results[$row]
    for ($i = 0; $i -lt $rowValues.Length; $i++) {
        $col = 2 + $i   # Patient1 starts at column B (2)
        $ws.Cells.Item($row, $col).Value = $rowValues[$i]
    }
}
